$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.202.14"
$ws.Range("E2").Value = "  -5.90%  "
$ws.Range("D3").Value = "2.217.23"
$ws.Range("E3").Value = "  -5.84%  "
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.10"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.619"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -7.37%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "69.88"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -5.78%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.555"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -7.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.70"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.61%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0950"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -7.09%  "
$ws.Range("E12").Value = "  -2.49%  "
$ws.Range("E13").Value = "  -3.72%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.77"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -7.30%  "
$ws.Range("D15").Value = "2.547.50"
$ws.Range("E15").Value = "  -6.08%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.81"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -9.58%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.842"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -8.84%  "
$ws.Range("D18").Value = "2.217.08"
$ws.Range("E18").Value = "  -6.18%  "
$ws.Range("D19").Value = "41.270.61"
$ws.Range("E19").Value = "  -5.72%  "
$ws.Range("E20").Value = "  -7.96%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.54"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.07"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -7.89%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "231.62"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -8.57%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.09"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +12.14%  "
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.70"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.49%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.42"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.84%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.85"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -7.38%  "
$ws.Range("E29").Value = "  -5.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "171.67"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.88%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.51"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -7.98%  "
$ws.Range("E32").Value = "  -7.89%  "
$ws.Range("E33").Value = "  -6.95%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0712"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.57%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.24"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.64%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.60"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -10.17%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.89"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.59%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "24.07"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +15.90%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0275"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.27%  "
$ws.Range("E40").Value = "  -5.37%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.84"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -11.44%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "65.64"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.21%  "
$ws.Range("B43").Value = "FTXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.02"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -10.04%  "
$ws.Range("B44").Value = "Algorand"
$ws.Range("C44").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.205"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.41%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.88"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.92%  "
$ws.Range("B46").Value = "Celestia"
$ws.Range("C46").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.89"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +11.56%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.100"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.47%  "
$ws.Range("B48").Value = "BinanceUSD"
$ws.Range("C48").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.00"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.09%  "
$ws.Range("B49").Value = "SynthetixNetwork"
$ws.Range("C49").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.57"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.96%  "
$ws.Range("E50").Value = "  -6.12%  "
$ws.Range("B51").Value = "ARBITRUM"
$ws.Range("C51").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.09"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.40%  "
